# Add new weekly price records for "Hortaliza, Terminal La Palmera de La Serena - Brocoli".
# Two new rows (Primera / Segunda quality grades) are inserted at the top of the
# existing date-ordered data block (before row 251), pushing all subsequent rows
# down by two, and the sheet grows from 378 to 380 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 251, shifting rows 251-378 down
# to 253-380. Excel copies the formatting of the row above, which already carries
# the date-style (s="2") on column D, matching the rest of the table.
$ws.Rows.Item(251).Insert()
$ws.Rows.Item(252).Insert()

# Row 251: "Primera" quality record for 2021-09-16 (serial 44455)
$ws.Range("A251").Value = 8
$ws.Range("B251").Value = "Terminal La Palmera de La Serena"
$ws.Range("C251").Value = "Coquimbo"
$ws.Range("D251").Value = 44455
$ws.Range("E251").Value = 4
$ws.Range("F251").Value = 100112023
$ws.Range("G251").Value = "Brócoli"
$ws.Range("H251").Value = "Sin especificar"
$ws.Range("I251").Value = "Primera"
$ws.Range("J251").Value = 2200
$ws.Range("K251").Value = 600
$ws.Range("L251").Value = 700
$ws.Range("M251").Value = 650
$ws.Range("N251").Value = "`$/unidad"
$ws.Range("O251").Value = "Provincia del Elquí"
$ws.Range("P251").Value = 650
$ws.Range("Q251").Value = 1
$ws.Range("R251").Value = "Hortaliza"

# Row 252: "Segunda" quality record for 2021-09-16 (serial 44455)
$ws.Range("A252").Value = 8
$ws.Range("B252").Value = "Terminal La Palmera de La Serena"
$ws.Range("C252").Value = "Coquimbo"
$ws.Range("D252").Value = 44455
$ws.Range("E252").Value = 4
$ws.Range("F252").Value = 100112023
$ws.Range("G252").Value = "Brócoli"
$ws.Range("H252").Value = "Sin especificar"
$ws.Range("I252").Value = "Segunda"
$ws.Range("J252").Value = 1400
$ws.Range("K252").Value = 500
$ws.Range("L252").Value = 550
$ws.Range("M252").Value = 525
$ws.Range("N252").Value = "`$/unidad"
$ws.Range("O252").Value = "Provincia del Elquí"
$ws.Range("P252").Value = 525
$ws.Range("Q252").Value = 1
$ws.Range("R252").Value = "Hortaliza"
